# Insert a new daily-revenue record for day 14 of July/2025 into the
# faturamento_diario sheet. The new row is inserted at worksheet row 15
# (right after the last existing July/2025 entry, day 13), which pushes
# every subsequent row down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at position 15, shifting rows 15..105 down to 16..106.
$ws.Rows.Item(15).Insert()

# Populate the newly inserted row with the new record.
$ws.Range("A15").Value2 = 14
$ws.Range("B15").Value2 = 21278.41
$ws.Range("C15").Value2 = 7
$ws.Range("D15").Value2 = 2025
$ws.Range("E15").Value2 = "07/2025"
